$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 currently only has C16, D16 (shared formula), E16 populated.
# Add the new sample data for row 16 (date 2018-04-02, serial 43192), matching
# the pattern used by the other rows (A col = date w/ date style, B col = CRM
# value, F/G cols = shared-string flags "With Junk" / "end of sample").

$ws.Range("A16").Value = 43192
$ws.Range("A16").NumberFormat = "m/d/yy"

$ws.Range("B16").Value = 2222.9530294207898

$ws.Range("F16").Value = "With Junk"
$ws.Range("G16").Value = "end of sample"

$wb.Save()
